# Angular material and bootstrap
# Inserts a new "Sheet2" worksheet (project/team roster) between "Team 1" and
# "Team 2", and updates a handful of cached selections / cell edits across
# the existing "Team *" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet right after "Team 1" (Excel names it "Sheet2"
#    automatically since "Sheet1" is already taken).
# ---------------------------------------------------------------------------
$team1 = $wb.Worksheets.Item("Team 1")
$newSheet = $wb.Worksheets.Add($null, $team1)

$newSheet.Range("A1").Value = "Team A"
$newSheet.Range("B1").Value = "Typeing Project"
$newSheet.Range("A2").Value = "Team B"
$newSheet.Range("B2").Value = "Trip Planner APP"
$newSheet.Range("A3").Value = "Team C"
$newSheet.Range("B3").Value = "Donation APP"
$newSheet.Range("A4").Value = "Team D"
$newSheet.Range("A5").Value = "Team E"
$newSheet.Range("C1").Value = "TypingDB"

$newSheet.Columns.Item(2).ColumnWidth = 31.5
$newSheet.Columns.Item(3).ColumnWidth = 24.166666666666668

$newSheet.Range("B7").Select()

# ---------------------------------------------------------------------------
# 2. "Team 1" - only the cached selection changes (D18 -> H9).
# ---------------------------------------------------------------------------
$team1.Range("H9").Select()

# ---------------------------------------------------------------------------
# 3. "Team 2" - cached selection changes (D11 -> F5).
# ---------------------------------------------------------------------------
$team2 = $wb.Worksheets.Item("Team 2")
$team2.Range("F5").Select()

# ---------------------------------------------------------------------------
# 4. "Team 3" - cached selection changes (C2 -> A5) and the A5 cell (which
#    held the stray "JT" label) is cleared, keeping its formatting.
# ---------------------------------------------------------------------------
$team3 = $wb.Worksheets.Item("Team 3")
$team3.Range("A5").ClearContents()
$team3.Range("A5").Select()

# ---------------------------------------------------------------------------
# 5. "Team 4" - cached selection changes (B14 -> C11) and A3 now holds the
#    "JT" label (previously "Phillip Benoit") using the same border/format
#    as the "JT" cell used to have on "Team 3" (no top/bottom border).
# ---------------------------------------------------------------------------
$team4 = $wb.Worksheets.Item("Team 4")
$team3.Range("A4").Copy()
$team4.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$team4.Range("A3").Value = "JT"
$team4.Range("C11").Select()

# ---------------------------------------------------------------------------
# 6. "Team 5" - cached selection changes (D13 -> ctrl-click union of
#    A1:A4, active cell A3).
# ---------------------------------------------------------------------------
$team5 = $wb.Worksheets.Item("Team 5")
$team5.Range("A3").Select()

# Restore "Team 1" as the active sheet/tab, matching the saved view state.
$team1.Activate()
$team1.Range("H9").Select()
